$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per diff
$ws.Range("C2").Value = 0.8981853442305029
$ws.Range("C3").Value = 0.8317615358706635
$ws.Range("C4").Value = 0.829290758366595
$ws.Range("C5").Value = 0.8738570379763733
$ws.Range("C6").Value = 0.8044522445963421
$ws.Range("C7").Value = 0.4849404622927854
$ws.Range("C8").Value = 0.5406942701798411
$ws.Range("C9").Value = 0.3516709972224711
$ws.Range("C10").Value = 0.3482540594913062
$ws.Range("C11").Value = 0.7453400503778338
$ws.Range("C12").Value = 0.7825317575680971
$ws.Range("C13").Value = 0.6565192418150488
$ws.Range("C14").Value = 0.5504124656278644
$ws.Range("D14").Value = 0.565853486568427
$ws.Range("E14").Value = 0.5163223577522386
$ws.Range("F14").Value = 0.556581823309596
$ws.Range("G14").Value = 0.5499189170133258
$ws.Range("H14").Value = 0.467531551857858
$ws.Range("I14").Value = 0.484206444334767
$ws.Range("L14").Value = 0.5435380384967919
$ws.Range("M14").Value = 0.5782979623492914
$ws.Range("N14").Value = 0.571494042163153
$ws.Range("O14").Value = 0.9524430656419657
$ws.Range("P14").Value = 0.4936896284284002
$ws.Range("Q14").Value = 0.5566170767820631
$ws.Range("R14").Value = 0.4608686455615878
$ws.Range("S14").Value = 0.5665585560177677
$ws.Range("T14").Value = 0.483289854050624
$ws.Range("U14").Value = 0.4782486074878375
$ws.Range("C15").Value = 0.4023289774216497
$ws.Range("C16").Value = 0.605129705480788
$ws.Range("D16").Value = 0.655061439438268
$ws.Range("E16").Value = 0.6150770431051297
$ws.Range("F16").Value = 0.6090306221962161
$ws.Range("G16").Value = 0.6302906182952994
$ws.Range("H16").Value = 0.6277550224302711
$ws.Range("I16").Value = 0.6203432806709577
$ws.Range("L16").Value = 0.6715428125609518
$ws.Range("M16").Value = 0.6914374878096352
$ws.Range("N16").Value = 0.6979715233079774
$ws.Range("O16").Value = 0.5769455822118198
$ws.Range("P16").Value = 0.6303881412131851
$ws.Range("Q16").Value = 0.941876340940121
$ws.Range("R16").Value = 0.5473961380924517
$ws.Range("S16").Value = 0.6609128145114102
$ws.Range("T16").Value = 0.6403354788375268
$ws.Range("U16").Value = 0.5786034718158767
$ws.Range("C17").Value = 0.6261357378854625
$ws.Range("D17").Value = 0.6510187224669604
$ws.Range("E17").Value = 0.6603111233480177
$ws.Range("F17").Value = 0.6542194383259912
$ws.Range("G17").Value = 0.662961178414097
$ws.Range("H17").Value = 0.4665129405286343
$ws.Range("I17").Value = 0.5424008810572687
$ws.Range("L17").Value = 0.6717373348017621
$ws.Range("M17").Value = 0.5567524779735683
$ws.Range("N17").Value = 0.5299765969162996
$ws.Range("O17").Value = 0.4578400330396476
$ws.Range("P17").Value = 0.4740501101321586
$ws.Range("Q17").Value = 0.5228180066079295
$ws.Range("R17").Value = 0.9678207599118943
$ws.Range("S17").Value = 0.7954983480176211
$ws.Range("T17").Value = 0.6192868942731278
$ws.Range("U17").Value = 0.5195140418502202
$ws.Range("C18").Value = 0.559291801545881
$ws.Range("C19").Value = 0.2118071203244705
$ws.Range("C20").Value = 0.1608416869863632
$ws.Range("C22").Value = 0.5758477728349909
$ws.Range("D22").Value = 0.632858063060262
$ws.Range("E22").Value = 0.5644293865401874
$ws.Range("F22").Value = 0.5881539320583904
$ws.Range("G22").Value = 0.5922830668267715
$ws.Range("H22").Value = 0.5283393504809212
$ws.Range("I22").Value = 0.5143876306739299
$ws.Range("L22").Value = 0.6188900186311833
$ws.Range("M22").Value = 0.5764117621486423
$ws.Range("N22").Value = 0.5473482185010652
$ws.Range("O22").Value = 0.5049256073307429
$ws.Range("P22").Value = 0.5228564495970084
$ws.Range("Q22").Value = 0.543724842078163
$ws.Range("R22").Value = 0.5263636166979903
$ws.Range("S22").Value = 0.6266189199139993
$ws.Range("T22").Value = 0.5782573269068908
$ws.Range("U22").Value = 0.5119837581136053
$ws.Range("C27").Value = 0.8348403942024935
$ws.Range("C28").Value = 0.4941799679128475
$ws.Range("C29").Value = 0.5993844295826897
$ws.Range("D29").Value = 0.6482120327705737
$ws.Range("E29").Value = 0.5782686598708006
$ws.Range("F29").Value = 0.6162166904298669
$ws.Range("G29").Value = 0.5958409978948308
$ws.Range("H29").Value = 0.5665753440564864
$ws.Range("I29").Value = 0.5338326257046987
$ws.Range("L29").Value = 0.6287557182982682
$ws.Range("M29").Value = 0.6145380010075698
$ws.Range("N29").Value = 0.635540868635675
$ws.Range("O29").Value = 0.5832084255544703
$ws.Range("P29").Value = 0.623852733516761
$ws.Range("Q29").Value = 0.6353402473796275
$ws.Range("R29").Value = 0.5217817393025161
$ws.Range("S29").Value = 0.6422444759059244
$ws.Range("T29").Value = 0.5755394370979757
$ws.Range("U29").Value = 0.5300334843857233
$ws.Range("C30").Value = 0.3895190866855442
$ws.Range("D30").Value = 0.4921832274466594
$ws.Range("E30").Value = 0.4441006500182716
$ws.Range("F30").Value = 0.3791518308681749
$ws.Range("G30").Value = 0.4550156947950562
$ws.Range("H30").Value = 0.3399163921772965
$ws.Range("I30").Value = 0.4144526211212216
$ws.Range("L30").Value = 0.4959586316735424
$ws.Range("M30").Value = 0.4193255353394157
$ws.Range("N30").Value = 0.3968122138965062
$ws.Range("O30").Value = 0.3674983299302888
$ws.Range("P30").Value = 0.3811459571415703
$ws.Range("Q30").Value = 0.3893191024077015
$ws.Range("S30").Value = 0.4948288886083743
$ws.Range("T30").Value = 0.5280562230190872
$ws.Range("U30").Value = 0.4695403205849761
$ws.Range("C32").Value = 0.5794809695958938
$ws.Range("D32").Value = 0.6361041034964985
$ws.Range("E32").Value = 0.5676838697943927
$ws.Range("F32").Value = 0.5896213976851233
$ws.Range("G32").Value = 0.595810366211031
$ws.Range("H32").Value = 0.5220537076330862
$ws.Range("I32").Value = 0.5101298184129563
$ws.Range("L32").Value = 0.6101061406290244
$ws.Range("M32").Value = 0.5704637627511275
$ws.Range("N32").Value = 0.5442310507775676
$ws.Range("O32").Value = 0.5019683843719066
$ws.Range("P32").Value = 0.5208207391702242
$ws.Range("Q32").Value = 0.5405853253597132
$ws.Range("R32").Value = 0.5110422343518973
$ws.Range("S32").Value = 0.6134399167834369
$ws.Range("T32").Value = 0.5732372165181104
$ws.Range("U32").Value = 0.5077394143607424
$ws.Range("C37").Value = 0.8099487324148484
$ws.Range("C38").Value = 0.5254331729390094
$ws.Range("C39").Value = 0.6075845621097523
$ws.Range("D39").Value = 0.5763878960198178
$ws.Range("E39").Value = 0.6184960552188207
$ws.Range("F39").Value = 0.5673997841730349
$ws.Range("C40").Value = 0.4319940979627413
$ws.Range("D40").Value = 0.4167758816573535
$ws.Range("E40").Value = 0.3908202277430965
$ws.Range("F40").Value = 0.4663938063623745
$ws.Range("C42").Value = 0.5216706110038343
$ws.Range("D42").Value = 0.5425332193688454
$ws.Range("E42").Value = 0.5079864515752036
$ws.Range("F42").Value = 0.5796883252172709

# Style swaps: move the bold+underline "row max" highlight to new max column
# Row 14: M14 (was max) -> normal; O14 (new max) -> bold+underline
$ws.Range("M14").Font.Bold = $false
$ws.Range("M14").Font.Underline = -4142
$ws.Range("O14").Font.Bold = $true
$ws.Range("O14").Font.Underline = 2

# Row 17: L17 (was max) -> normal; R17 (new max) -> bold+underline
$ws.Range("L17").Font.Bold = $false
$ws.Range("L17").Font.Underline = -4142
$ws.Range("R17").Font.Bold = $true
$ws.Range("R17").Font.Underline = 2
